$d = $word.ActiveDocument

# 1) "类型:普通" -> "类型:红包", with rFonts hint=eastAsia on the new run
$r1 = $d.Content
$r1.Find.Execute("普通", $true, $false, $false, $false, $false, $true, 1, $false, "红包", 2) | Out-Null

# 2) "照片:红包再猜猜看" -> "照片:再猜猜看！" (remove 红包 from first run, append ！ as new run)
$r2 = $d.Content
$r2.Find.Execute("照片:红包", $true, $false, $false, $false, $false, $true, 1, $false, "照片:", 2) | Out-Null

$r3 = $d.Content
$r3.Find.Execute("再猜猜看", $true, $false, $false, $false, $false, $true, 1, $false, "再猜猜看！", 2) | Out-Null
